$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.354.52'
$ws.Cells.Item(2, 5).Value = '  +0.52%  '

$ws.Cells.Item(3, 4).Value = '1.592.02'

$ws.Cells.Item(4, 5).Value = '  -0.34%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '211.40'

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.506'
$ws.Cells.Item(6, 5).Value = '  +0.10%  '

$ws.Cells.Item(7, 5).Value = '  -0.33%  '

$ws.Cells.Item(8, 5).Value = '  +0.06%  '

$ws.Cells.Item(9, 5).Value = '  -0.17%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '19.47'
$ws.Cells.Item(10, 5).Value = '  -0.49%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0846'
$ws.Cells.Item(11, 5).Value = '  +0.03%  '

$ws.Cells.Item(12, 4).Value = '1.815.57'
$ws.Cells.Item(12, 5).Value = '  +0.55%  '

$ws.Cells.Item(13, 4).Value = '1.615.34'
$ws.Cells.Item(13, 5).Value = '  +2.22%  '

$ws.Cells.Item(14, 5).Value = '  +0.88%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.522'
$ws.Cells.Item(15, 5).Value = '  +0.97%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '64.67'
$ws.Cells.Item(16, 5).Value = '  -0.06%  '

$ws.Cells.Item(17, 4).Value = '26.355.14'
$ws.Cells.Item(17, 5).Value = '  +0.42%  '

$ws.Cells.Item(18, 4).Value = '0.0₃0732'
$ws.Cells.Item(18, 5).Value = '  -0.99%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '7.48'
$ws.Cells.Item(19, 5).Value = '  +3.78%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '212.13'
$ws.Cells.Item(20, 5).Value = '  +2.63%  '

$ws.Cells.Item(21, 5).Value = '  -0.29%  '

$ws.Cells.Item(22, 5).Value = '  +0.93%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '2.19'
$ws.Cells.Item(23, 5).Value = '  -1.01%  '

$ws.Cells.Item(24, 5).Value = '  +1.86%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '144.86'
$ws.Cells.Item(25, 5).Value = '  +0.34%  '

$ws.Cells.Item(26, 5).Value = '  -0.36%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '7.06'
$ws.Cells.Item(27, 5).Value = '  +0.67%  '

$ws.Cells.Item(28, 5).Value = '  -0.65%  '

$ws.Cells.Item(29, 5).Value = '  -0.30%  '

$ws.Cells.Item(30, 5).Value = '  -0.19%  '

$ws.Cells.Item(31, 5).Value = '  +0.87%  '

$ws.Cells.Item(32, 5).Value = '  -0.39%  '

$ws.Cells.Item(33, 5).Value = '  +1.09%  '

$ws.Cells.Item(34, 4).Value = '1.342.06'
$ws.Cells.Item(34, 5).Value = '  +4.42%  '

$ws.Cells.Item(35, 5).Value = '  -1.19%  '

$ws.Cells.Item(36, 5).Value = '  -1.29%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.0167'
$ws.Cells.Item(38, 5).Value = '  -0.06%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.06'
$ws.Cells.Item(39, 5).Value = '  -15.02%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.817'
$ws.Cells.Item(40, 5).Value = '  +0.19%  '

$ws.Cells.Item(41, 5).Value = '  +4.03%  '

$ws.Cells.Item(42, 5).Value = '  -0.26%  '

$ws.Cells.Item(43, 5).Value = '  +0.11%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.763'
$ws.Cells.Item(44, 5).Value = '  -0.49%  '

$ws.Cells.Item(45, 4).Value = '1.728.61'
$ws.Cells.Item(45, 5).Value = '  +0.54%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '61.63'
$ws.Cells.Item(46, 5).Value = '  -0.93%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '88.05'
$ws.Cells.Item(47, 5).Value = '  -1.03%  '

$ws.Cells.Item(48, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(48, 4).Value = '0.0₆0105'
$ws.Cells.Item(48, 5).Value = '  +2.17%  '

$ws.Cells.Item(49, 2).Value = 'RenderToken'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.51'
$ws.Cells.Item(49, 5).Value = '  -2.98%  '

$ws.Cells.Item(50, 2).Value = 'Algorand'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.0985'
$ws.Cells.Item(50, 5).Value = '  -3.40%  '

$ws.Cells.Item(51, 2).Value = 'Cronos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.0505'
$ws.Cells.Item(51, 5).Value = '  -1.06%  '
